$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F column) counts for sheets "展览" and "全部类型"
$updates = @{
    3  = 84
    4  = 255
    5  = 41
    6  = 519
    7  = 47
    8  = 1967
    11 = 4155
    13 = 270
    16 = 15
    20 = 402
    25 = 57
    26 = 6
    29 = 188
    30 = 295
    31 = 1631
    32 = 233
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# F18 differs slightly between the two sheets
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F18").Value = 2872

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F18").Value = 2873
